$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 273, shifting existing rows 273:325 down to 274:326
$ws.Rows(273).Insert()

# Populate the newly inserted row 273 with the new weekly record
$ws.Range("A273").Value = 7
$ws.Range("B273").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C273").Value = "Ñuble"
$ws.Range("D273").Value = 44522
$ws.Range("E273").Value = 16
$ws.Range("F273").Value = 100112020
$ws.Range("G273").Value = "Tomate"
$ws.Range("H273").Value = "Larga vida"
$ws.Range("I273").Value = "Primera"
$ws.Range("J273").Value = 360
$ws.Range("K273").Value = 7000
$ws.Range("L273").Value = 8000
$ws.Range("M273").Value = 7500
$ws.Range("N273").Value = "$/caja 15 kilos"
$ws.Range("O273").Value = "Región del Maule"
$ws.Range("P273").Value = 500
$ws.Range("Q273").Value = 15
$ws.Range("R273").Value = "Hortaliza"
